# Add a "shortList" boolean column (E) to the population sheet, flagging a
# shortlist of countries for future work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, matching the style already used by the other
# header cells in row 1.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "shortList"

# Copy the standard data-row style (borders/centering) down to the new
# column E before filling in the boolean values. Use the plain (non
# highlighted) style shared by most rows -- D2 -- as the source, rather than
# the per-row D column (a couple of rows, e.g. the "home" countries, use a
# highlighted style that the new column does not pick up).
$ws.Range("D2").Copy()
$ws.Range("E2:E87").PasteSpecial(-4122)

$shortList = @($false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$true,$true,$false,$false,$false,$false,$false,$false,$false,$true,$false,$true,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false,$false,$true,$false,$false,$false)

for ($i = 0; $i -lt $shortList.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $shortList[$i]
}

# Restore the selection/scroll state recorded for this edit.
[void]$ws.Range("A57").Select()
[void]$ws.Range("F64").Select()
